# Auto update Excel log
# Appends new sensor-log rows to five worksheets: PIR, Humidity, Temperature,
# Proximity and mmWave. Columns (A:F) = Date, Timestamp, Hour, Location,
# Value, Status.
#
# New row data is passed as a single FLAT value list (rather than an array
# of per-row arrays) because a one-row "array of arrays" gets silently
# unrolled by this host's PowerShell-ish interpreter, which would scramble
# a single appended row across cells/characters. Using one flat list with a
# known column stride sidesteps that ambiguity for both single- and
# multi-row appends.
#
# $TextCols lists the 0-based column indexes that must be forced to text
# number-format before the value is written, so Excel's autodetection
# doesn't silently turn a literal "2026-01-30" into a date serial or a
# literal "87.2%" into a numeric percentage -- the source log stores these
# as plain strings (inline/shared text), never as typed numbers/dates.

$wb = $excel.ActiveWorkbook

function Append-Rows {
    param(
        [string]$SheetName,
        [int]$StartRow,
        [int]$NumCols,
        [object[]]$Flat,
        [int[]]$TextCols
    )

    $ws = $wb.Worksheets.Item($SheetName)
    $numRows = $Flat.Count / $NumCols
    $r = $StartRow
    $i = 0
    for ($k = 0; $k -lt $numRows; $k = $k + 1) {
        for ($c = 0; $c -lt $NumCols; $c = $c + 1) {
            $cell = $ws.Cells.Item($r, $c + 1)
            if ($TextCols -contains $c) {
                $cell.NumberFormat = "@"
            }
            $cell.Value = $Flat[$i]
            $i = $i + 1
        }
        $r = $r + 1
    }
}

# ---------------------------------------------------------------------
# PIR sheet: rows 114-125 (Bathroom / No Motion / Inactive)
# ---------------------------------------------------------------------
$pirFlat = @(
    "2026-01-30", "13:08:50", "13:00", "Bathroom", "No Motion", "Inactive",
    "2026-01-30", "13:08:51", "13:00", "Bathroom", "No Motion", "Inactive",
    "2026-01-30", "13:08:55", "13:00", "Bathroom", "No Motion", "Inactive",
    "2026-01-30", "13:09:00", "13:00", "Bathroom", "No Motion", "Inactive",
    "2026-01-30", "13:09:05", "13:00", "Bathroom", "No Motion", "Inactive",
    "2026-01-30", "13:09:10", "13:00", "Bathroom", "No Motion", "Inactive",
    "2026-01-30", "13:09:15", "13:00", "Bathroom", "No Motion", "Inactive",
    "2026-01-30", "13:09:20", "13:00", "Bathroom", "No Motion", "Inactive",
    "2026-01-30", "13:09:25", "13:00", "Bathroom", "No Motion", "Inactive",
    "2026-01-30", "13:09:30", "13:00", "Bathroom", "No Motion", "Inactive",
    "2026-01-30", "13:09:35", "13:00", "Bathroom", "No Motion", "Inactive",
    "2026-01-30", "13:09:40", "13:00", "Bathroom", "No Motion", "Inactive"
)
Append-Rows "PIR" 114 6 $pirFlat @(0)

# ---------------------------------------------------------------------
# Humidity sheet: row 57
# ---------------------------------------------------------------------
$humidityFlat = @(
    "2026-01-30", "13:09:07", "13:00", "Bathroom", "87.2%", "Active"
)
Append-Rows "Humidity" 57 6 $humidityFlat @(0, 4)

# ---------------------------------------------------------------------
# Temperature sheet: row 57
# ---------------------------------------------------------------------
$temperatureFlat = @(
    "2026-01-30", "13:09:07", "13:00", "Bathroom", "22.6C", "Active"
)
Append-Rows "Temperature" 57 6 $temperatureFlat @(0)

# ---------------------------------------------------------------------
# Proximity sheet: rows 49-52 (Bathroom Door ENTER/EXIT events)
# ---------------------------------------------------------------------
$proximityFlat = @(
    "2026-01-30", "13:08:50", "13:00", "Bathroom Door", "ENTER", "User ENTERED Bathroom",
    "2026-01-30", "13:08:56", "13:00", "Bathroom Door", "EXIT", "User EXITED Bathroom",
    "2026-01-30", "13:09:24", "13:00", "Bathroom Door", "ENTER", "User ENTERED Bathroom",
    "2026-01-30", "13:09:39", "13:00", "Bathroom Door", "EXIT", "User EXITED Bathroom"
)
Append-Rows "Proximity" 49 6 $proximityFlat @(0)

# ---------------------------------------------------------------------
# mmWave sheet: rows 41-46 (Living Room emergency/motion events)
# ---------------------------------------------------------------------
$mmWaveFlat = @(
    "2026-01-30", "14:03:19", "14:00", "Living Room", "FALL_DETECTED", "EMERGENCY",
    "2026-01-30", "14:03:19", "14:00", "Living Room", "FALL_DETECTED", "EMERGENCY",
    "2026-01-30", "14:03:31", "14:00", "Living Room", "NO_MOTION_DETECTED", "Inactive",
    "2026-01-30", "14:03:42", "14:00", "Living Room", "PRESENCE_DETECTED", "Active",
    "2026-01-30", "14:03:52", "14:00", "Living Room", "PRESENCE_DETECTED", "Active",
    "2026-01-30", "14:04:09", "14:00", "Living Room", "PRESENCE_DETECTED", "Active"
)
Append-Rows "mmWave" 41 6 $mmWaveFlat @(0)
